$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.103.66"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "'2.587.17"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'527.65"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'138.96"
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "'2.598.71"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "'6.42"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("D14").Value = "'3.044.45"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "'59.073.20"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "'20.48"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'2.586.78"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "'343.47"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "'10.06"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'6.41"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'66.73"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'7.06"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'0.0₃0719"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "'5.89"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").Value = "'18.68"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'149.46"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'3.96"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "'36.80"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'0.824"
$ws.Range("E39").Value = "  -5.14%  "
$ws.Range("D40").Value = "'0.810"
$ws.Range("E40").Value = "  -6.72%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'0.600"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'10.78"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'268.23"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'0.0952"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "'0.0513"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "'1.959.51"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'18.17"
$ws.Range("E51").Value = "  -2.94%  "
